# Auto-generated edit script: apply scheduled market-price refresh values
# to the Maduin_Profits workbook (columns H-N per Leve row) across all 8
# job-class sheets, matching the upstream runner's commit.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 306
$ws.Range("I9").Value = 74.666664
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 74.666664
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = 94.333336
$ws.Range("N9").Value = -1338
$ws.Range("H51").Value = 11562.125
$ws.Range("I51").Value = 9499.6
$ws.Range("K51").Value = 9499.6
$ws.Range("M51").Value = -9015.6
$ws.Range("H127").Value = 4999
$ws.Range("I127").Value = 4999
$ws.Range("K127").Value = 14997
$ws.Range("M127").Value = -10037

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1549.9474
$ws.Range("I97").Value = 1450.8
$ws.Range("K97").Value = 1450.8
$ws.Range("M97").Value = -954.8

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 20045.8
$ws.Range("I11").Value = 150
$ws.Range("J11").Value = 25019.75
$ws.Range("K11").Value = 150
$ws.Range("L11").Value = 25019.75
$ws.Range("M11").Value = -10
$ws.Range("N11").Value = -25299.75
$ws.Range("H37").Value = 2504.2
$ws.Range("I37").Value = 1880.25
$ws.Range("J37").Value = 5000
$ws.Range("K37").Value = 1880.25
$ws.Range("L37").Value = 5000
$ws.Range("M37").Value = -1743.25
$ws.Range("N37").Value = -5274
$ws.Range("H94").Value = 2886.8572
$ws.Range("I94").Value = 2041.8
$ws.Range("J94").Value = 4999.5
$ws.Range("K94").Value = 2041.8
$ws.Range("L94").Value = 4999.5
$ws.Range("M94").Value = -1590.8
$ws.Range("N94").Value = -5901.5
$ws.Range("H99").Value = 1740
$ws.Range("I99").Value = 1740
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1740
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -242
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 2249.25
$ws.Range("I105").Value = 2249.25
$ws.Range("K105").Value = 2249.25
$ws.Range("M105").Value = -502.25

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 733.1818
$ws.Range("I22").Value = 607
$ws.Range("K22").Value = 607
$ws.Range("M22").Value = -257
$ws.Range("H92").Value = 57989
$ws.Range("J92").Value = 57989
$ws.Range("L92").Value = 57989
$ws.Range("N92").Value = -62981

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 57928.473
$ws.Range("I2").Value = 21.1875
$ws.Range("J2").Value = 366767.34
$ws.Range("K2").Value = 127.125
$ws.Range("L2").Value = 2200604.04
$ws.Range("M2").Value = -14.125
$ws.Range("N2").Value = -2200830.04
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H11").Value = 15000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 15000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 45000
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -45280
$ws.Range("H70").Value = 3900
$ws.Range("I70").Value = 3900
$ws.Range("K70").Value = 11700
$ws.Range("M70").Value = -11385
$ws.Range("H73").Value = 3900
$ws.Range("I73").Value = 3900
$ws.Range("K73").Value = 11700
$ws.Range("M73").Value = -10608
$ws.Range("H80").Value = 2599.5
$ws.Range("I80").Value = 2996
$ws.Range("J80").Value = 2203
$ws.Range("K80").Value = 8988
$ws.Range("L80").Value = 6609
$ws.Range("M80").Value = -8052
$ws.Range("N80").Value = -8481
$ws.Range("H83").Value = 2599.5
$ws.Range("I83").Value = 2996
$ws.Range("J83").Value = 2203
$ws.Range("K83").Value = 26964
$ws.Range("L83").Value = 19827
$ws.Range("M83").Value = -22284
$ws.Range("N83").Value = -29187
$ws.Range("H92").Value = 657
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 774.75
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 2324.25
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -4820.25
$ws.Range("H99").Value = 25
$ws.Range("I99").Value = 25
$ws.Range("K99").Value = 75
$ws.Range("M99").Value = 2171
$ws.Range("H107").Value = 1951.5
$ws.Range("I107").Value = 1803
$ws.Range("J107").Value = 2100
$ws.Range("K107").Value = 5409
$ws.Range("L107").Value = 6300
$ws.Range("M107").Value = -3489
$ws.Range("N107").Value = -10140
$ws.Range("H131").Value = 838
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H137").Value = 3993.7144
$ws.Range("I137").Value = 3999.5
$ws.Range("J137").Value = 3991.4
$ws.Range("K137").Value = 11998.5
$ws.Range("L137").Value = 11974.2
$ws.Range("M137").Value = -6898.5
$ws.Range("N137").Value = -22174.2

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8347.5
$ws.Range("I80").Value = 7686
$ws.Range("J80").Value = 10332
$ws.Range("K80").Value = 7686
$ws.Range("L80").Value = 10332
$ws.Range("M80").Value = -6688
$ws.Range("N80").Value = -12328
$ws.Range("H83").Value = 8347.5
$ws.Range("I83").Value = 7686
$ws.Range("J83").Value = 10332
$ws.Range("K83").Value = 38430
$ws.Range("L83").Value = 51660
$ws.Range("M83").Value = -33438
$ws.Range("N83").Value = -61644
$ws.Range("H102").Value = 1318.1
$ws.Range("I102").Value = 966.7143
$ws.Range("J102").Value = 2138
$ws.Range("K102").Value = 966.7143
$ws.Range("L102").Value = 2138
$ws.Range("M102").Value = 655.2857
$ws.Range("N102").Value = -5382

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5110.1113
$ws.Range("I22").Value = 3748.6667
$ws.Range("J22").Value = 7833
$ws.Range("K22").Value = 3748.6667
$ws.Range("L22").Value = 7833
$ws.Range("M22").Value = -3453.6667
$ws.Range("N22").Value = -8423
$ws.Range("H27").Value = 5110.1113
$ws.Range("I27").Value = 3748.6667
$ws.Range("J27").Value = 7833
$ws.Range("K27").Value = 3748.6667
$ws.Range("L27").Value = 7833
$ws.Range("M27").Value = -3641.6667
$ws.Range("N27").Value = -8047

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7000
$ws.Range("I96").Value = 8334
$ws.Range("J96").Value = 2998
$ws.Range("K96").Value = 8334
$ws.Range("L96").Value = 2998
$ws.Range("M96").Value = -6961
$ws.Range("N96").Value = -5744
$ws.Range("H122").Value = 1568.8572
$ws.Range("I122").Value = 1568.8572
$ws.Range("K122").Value = 4706.571599999999
$ws.Range("M122").Value = -2256.571599999999
$ws.Range("H126").Value = 827.375
$ws.Range("J126").Value = 295.5
$ws.Range("L126").Value = 886.5
$ws.Range("N126").Value = -5826.5

